$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.098162608357640124,
    -0.0059999999675817151,
    -0.0039999999747166726,
    -0.0079999999525917076,
    -0.0029999999785204068,
    -0.0019999999824520387,
    -0.0099999999400206541,
    -0.0099999999410451679,
    -0.0019999999860438322,
    -0.001999999988795409,
    0.040126609686510051,
    -0.0034999999808884574,
    -0.0034999999800104931,
    -0.0079999999559561275,
    -0.0009999999933887338,
    -0.0019999999884316999,
    0.011588820192413252,
    -0.003999999978631763,
    -0.003999999978103741,
    -0.0039999999763598026,
    -0.0039999999760960137,
    -0.003999999975943247,
    -0.0049999999676400009,
    -0.0199999998845799,
    -0.019999999882991837,
    0.033233986034707286,
    -0.0024999999787942428,
    -0.0019999999765936138,
    -0.0069999999465863993,
    -0.016255773729203948,
    -0.0069999999429466442,
    -0.0099999999266451312,
    -0.0039999999584807711
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
